# The "YIELD" column (Q) data was removed from the bond table: the header
# text in Q1 and the computed yield values in Q2:Q12 are cleared, leaving
# the cells (and their existing number-format styles) behind. The column Q
# is left selected afterwards, matching the saved selection state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("Q:Q").Select()
$ws.Range("Q1:Q12").ClearContents()
